# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new value
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value  = 262
$wsExhibition.Range("F12").Value = 109
$wsExhibition.Range("F13").Value = 2383
$wsExhibition.Range("F15").Value = 40
$wsExhibition.Range("F17").Value = 549
$wsExhibition.Range("F21").Value = 49
$wsExhibition.Range("F22").Value = 1889
$wsExhibition.Range("F23").Value = 4030
$wsExhibition.Range("F24").Value = 31
$wsExhibition.Range("F26").Value = 1185
$wsExhibition.Range("F28").Value = 2093
$wsExhibition.Range("F36").Value = 695
$wsExhibition.Range("F38").Value = 416

# Sheet "全部类型" (all types) - same events, shifted by one extra row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 262
$wsAll.Range("F12").Value = 109
$wsAll.Range("F13").Value = 2383
$wsAll.Range("F16").Value = 40
$wsAll.Range("F18").Value = 549
$wsAll.Range("F22").Value = 49
$wsAll.Range("F23").Value = 1889
$wsAll.Range("F24").Value = 4030
$wsAll.Range("F25").Value = 31
$wsAll.Range("F27").Value = 1185
$wsAll.Range("F29").Value = 2093
$wsAll.Range("F37").Value = 695
$wsAll.Range("F39").Value = 416
